$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the image-path strings (shared strings) referenced by A2 and A3
$ws.Range("A2").Value = "Rating\ratingCS+3.png"
$ws.Range("A3").Value = "Rating\ratingCS-3.png"

# Move the active selection from A3 to C7
$ws.Range("C7").Select() | Out-Null
